$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# Add new fuel-source rows mirroring existing priority values:
#   row 15 "crude oil"                   -> same priority as "petroleum" (row 11)
#   row 16 "heavy or residual fuel oil"  -> same priority as "petroleum" (row 11)
#   row 17 "municipal solid waste"       -> same priority as "biomass"   (row 9)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# Add a header label above the year columns and format it
$ws.Range("A1").Value = "Dispatch Priority (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Widen column A slightly to better fit the new header text
$ws.Columns.Item(1).ColumnWidth = 23
